# ccaa_t48 add f_s2 data.
# Extend the "cht12_f_s" sheet (x vs f_s lookup table) from 3 rows to 32
# rows, updating B2 and replacing the old A3/B3 pair with the full
# 31-point curve, then make this sheet the active tab/selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cht12_f_s")

$aVals = @(1, 1.0475651189127899, 1.0679501698754199, 1.0973952434881, 1.1472253680634199, 1.19705549263873, 1.2468856172140399, 1.29671574178935, 1.3488108720271801, 1.4009060022649999, 1.4484711211778001, 1.49830124575311, 1.55039637599094, 1.6002265005662499, 1.65005662514156, 1.6998867497168699, 1.74971687429218, 1.7995469988674899, 1.8810872027179999, 1.96715741789354, 2.0509626274065602, 2.14835787089467, 2.2321630804077, 2.3227633069082598, 2.4201585503963701, 2.5016987542468798, 2.6013590033975, 2.6761041902604701, 2.7848244620611502, 2.90033975084937, 3)
$bVals = @(0.91425269645608598, 0.92619414483821205, 0.93043143297380504, 0.93736517719568502, 0.94699537750385199, 0.95624036979969096, 0.96625577812018404, 0.97473035439137101, 0.98320493066255699, 0.99129429892141696, 0.99822804314329705, 1.0055469953775, 1.0120955315870499, 1.01941448382126, 1.0251926040061601, 1.03174114021571, 1.0371340523882799, 1.04291217257318, 1.05138674884437, 1.0602465331278801, 1.0675654853620899, 1.0760400616332799, 1.08258859784283, 1.0899075500770401, 1.09645608628659, 1.10184899845916, 1.10762711864406, 1.1118644067796599, 1.1172573189522299, 1.1226502311247999, 1.1268875192604)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "0.000"
    $ws.Cells.Item($row, 2).NumberFormat = "0.000"
}

$ws.Columns("A:B").ColumnWidth = 4.59

$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
